{"js": "// 1. Version \"4\" -> \"5\" (revision-history paragraph)\n{\n  const results = context.document.body.search(\"Version 4\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    const versionPara = results.items[0].paragraphs.getFirst();\n    const hit = versionPara.search(\"4\", { matchCase: true });\n    hit.load(\"items\");\n    await context.sync();\n    hit.items[0].insertText(\"5\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 2. Updated \"last modified\" date field text\n{\n  const results = context.document.body.search(\"11/15/21 9:55 AM\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"7/22/22 10:39 AM\", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// 3. \"ISO 21434\" -> bold \"ISO/SAE 21434\" (rest of sentence stays regular)\n{\n  const results = context.document.body.search(\"ISO 21434\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\"ISO/SAE 21434\", Word.InsertLocation.replace);\n    await context.sync();\n\n    const boldHit = context.document.body.search(\"ISO/SAE 21434\", { matchCase: true });\n    boldHit.load(\"items\");\n    await context.sync();\n    boldHit.items[0].font.bold = true;\n    await context.sync();\n  }\n}\n\n// 4. License text: \"(CC4-SA)\" -> \"(CC BY-SA-4.0)\" inside the existing bold run\n{\n  const results = context.document.body.search(\n    \"Creative Commons Attribution-Share Alike (CC4-SA)\",\n    { matchCase: true }\n  );\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\n      \"Creative Commons Attribution-Share Alike (CC BY-SA-4.0)\",\n      Word.InsertLocation.replace\n    );\n    await context.sync();\n  }\n}\n\n// 5. Append a new sentence (with a bold label + superscript reference) after the\n//    \"...be identified, updates will be made.\" sentence.\n{\n  const results = context.document.body.search(\n    \" be identified, updates will be made. \",\n    { matchCase: true }\n  );\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    results.items[0].insertText(\" be identified, updates will be made.\", Word.InsertLocation.replace);\n    await context.sync();\n\n    const anchorResults = context.document.body.search(\" be identified, updates will be made.\", { matchCase: true });\n    anchorResults.load(\"items\");\n    await context.sync();\n    const anchor = anchorResults.items[0];\n\n    const addedText =\n      \" The design security analysis should consider the best practices described in \" +\n      \"Secure Design Principles [2].\";\n    const addedRange = anchor.insertText(addedText, Word.InsertLocation.after);\n    await context.sync();\n\n    const boldLabel = addedRange.search(\"Secure Design Principles \", { matchCase: true });\n    boldLabel.load(\"items\");\n    await context.sync();\n    boldLabel.items[0].font.bold = true;\n    await context.sync();\n\n    const refMark = addedRange.search(\"[2]\", { matchCase: true });\n    refMark.load(\"items\");\n    await context.sync();\n    const refRange = refMark.items[0];\n    refRange.font.bold = true;\n    refRange.font.color = \"#0070C0\";\n    refRange.font.superscript = true;\n    await context.sync();\n  }\n}\n\n// 6. Add a new \"Secure Design Principles (AVCDL secondary document)\" item to the\n//    numbered References list, after the existing \"Security Requirements Taxonomy\" entry.\n{\n  const paragraphs = context.document.body.paragraphs;\n  paragraphs.load(\"items/text\");\n  await context.sync();\n\n  let referenceItem = null;\n  for (let i = 0; i < paragraphs.items.length; i++) {\n    if (paragraphs.items[i].text.indexOf(\"AVCDL secondary document\") !== -1) {\n      referenceItem = paragraphs.items[i];\n    }\n  }\n\n  if (referenceItem) {\n    const newPara = referenceItem.insertParagraph(\n      \"Secure Design Principles (AVCDL secondary document)\",\n      Word.InsertLocation.after\n    );\n    await context.sync();\n\n    const nonBold = newPara.search(\"(AVCDL secondary document)\", { matchCase: true });\n    nonBold.load(\"items\");\n    await context.sync();\n    nonBold.items[0].font.bold = false;\n    await context.sync();\n  }\n}\n", "ps1": "$d = $word.ActiveDocument\n$wdFindStop = 0\n$wdReplaceOne = 1\n$wdReplaceAll = 2\n\n# 1. Version \"4\" -> \"5\" (revision-history paragraph). Scope to that single\n#    paragraph so no other lone \"4\" in the document is touched.\n$verPara = $d.Paragraphs(3)\n$verRng = $verPara.Range\n$verRng.Find.Execute(\"4\", $false, $false, $false, $false, $false, $true, $wdFindStop, $false, \"5\", $wdReplaceAll)\n\n# 2. Updated \"last modified\" date field text\n$dateRng = $d.Content\n$dateRng.Find.Execute(\"11/15/21 9:55 AM\", $false, $false, $false, $false, $false, $true, $wdFindStop, $false, \"7/22/22 10:39 AM\", $wdReplaceAll)\n\n# 3. \"ISO 21434\" -> bold \"ISO/SAE 21434\" (rest of sentence stays regular)\n$isoRng = $d.Content\n$isoRng.Find.Execute(\"ISO 21434\", $false, $false, $false, $false, $false, $true, $wdFindStop, $false, \"ISO/SAE 21434\", $wdReplaceAll)\n\n$isoBoldRng = $d.Content\n$isoBoldRng.Find.Execute(\"ISO/SAE 21434\", $false, $false, $false, $false, $false, $true, $wdFindStop, $false, $null, 0)\n$isoBoldRng.Bold = $true\n\n# 4. License text: \"(CC4-SA)\" -> \"(CC BY-SA-4.0)\" inside the existing bold run\n$ccRng = $d.Content\n$ccRng.Find.Execute(\"Creative Commons Attribution-Share Alike (CC4-SA)\", $false, $false, $false, $false, $false, $true, $wdFindStop, $false, \"Creative Commons Attribution-Share Alike (CC BY-SA-4.0)\", $wdReplaceAll)\n\n# 5. Append a new sentence (with a bold label + superscript reference) after the\n#    \"...be identified, updates will be made.\" sentence.\n$trimRng = $d.Content\n$trimRng.Find.Execute(\" be identified, updates will be made. \", $false, $false, $false, $false, $false, $true, $wdFindStop, $false, \" be identified, updates will be made.\", $wdReplaceAll)\n\n$anchorRng = $d.Content\n$anchorRng.Find.Execute(\" be identified, updates will be made.\", $false, $false, $false, $false, $false, $true, $wdFindStop, $false, $null, 0)\n$anchorRng.Collapse(0)\n$anchorRng.InsertAfter(\" The design security analysis should consider the best practices described in Secure Design Principles [2].\")\n\n$boldLabelRng = $d.Content\n$boldLabelRng.Find.Execute(\"Secure Design Principles \", $false, $false, $false, $false, $false, $true, $wdFindStop, $false, $null, 0)\n$boldLabelRng.Bold = $true\n\n$refMarkRng = $d.Content\n$refMarkRng.Find.Execute(\"[2]\", $false, $false, $false, $false, $false, $true, $wdFindStop, $false, $null, 0)\n$refMarkRng.Bold = $true\n$refMarkRng.Font.Color = 0xC07000\n$refMarkRng.Font.Superscript = $true\n\n# 6. Add a new \"Secure Design Principles (AVCDL secondary document)\" item to the\n#    numbered References list, after the existing \"Security Requirements Taxonomy\" entry.\n$lastPara = $d.Paragraphs($d.Paragraphs.Count)\n$lastPara.Range.InsertParagraphAfter()\n$newPara = $d.Paragraphs($d.Paragraphs.Count)\n$newPara.Range.Text = \"Secure Design Principles (AVCDL secondary document)\"\n\n$newParaUnboldRng = $newPara.Range.Duplicate\n$newParaUnboldRng.Find.Execute(\"(AVCDL secondary document)\", $false, $false, $false, $false, $false, $true, $wdFindStop, $false, $null, 0)\n$newParaUnboldRng.Bold = $false\n"}
